$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: TC10 - wait / toWait
$ws.Cells.Item(11, 1).Value = "TC10"
$ws.Cells.Item(11, 2).Value = "wait"
$ws.Cells.Item(11, 3).Value = "toWait"

# Row 12: TC11 - Click on Logout (fill description/locator columns before the ID column,
# matching the shared-string discovery order recorded in the workbook)
$ws.Cells.Item(12, 2).Value = "Click on Logout"
$ws.Cells.Item(12, 4).Value = 'xpath~//*[@id="logoutLink"]'
$ws.Cells.Item(12, 3).Value = "click"
$ws.Cells.Item(12, 1).Value = "TC11"

# Row 13: TC12 - wait / toWait
$ws.Cells.Item(13, 1).Value = "TC12"
$ws.Cells.Item(13, 2).Value = "wait"
$ws.Cells.Item(13, 3).Value = "toWait"

# Update selection to match target state
$ws.Range("B13").Select()
